$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 233; everything below (233-243) shifts down to (234-244)
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with the new weekly record
$ws.Cells.Item(233, 1).Value = 3
$ws.Cells.Item(233, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(233, 3).Value = "Coquimbo"
$ws.Cells.Item(233, 4).Value = 45041
$ws.Cells.Item(233, 5).Value = 5
$ws.Cells.Item(233, 6).Value = 100112026
$ws.Cells.Item(233, 7).Value = "Haba"
$ws.Cells.Item(233, 8).Value = "Sin especificar"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 80
$ws.Cells.Item(233, 11).Value = 21000
$ws.Cells.Item(233, 12).Value = 22000
$ws.Cells.Item(233, 13).Value = 21500
$ws.Cells.Item(233, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(233, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(233, 16).Value = 860
$ws.Cells.Item(233, 17).Value = 25
$ws.Cells.Item(233, 18).Value = "Hortaliza"
